# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.360.82"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "3.619.10"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.213"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.646"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "4.194.71"
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "604.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "70.414.28"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.602.20"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").Value = "0.0₃0885"
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("D36").Value = "3.911.89"
$ws.Range("E36").Value = "  +5.18%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "519.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.76%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.51%  "
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.08%  "
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("E51").Value = "  +0.91%  "
